# "Generate Report for Archive"
# The localization status report is regenerated: every cell that showed
# "Ready for handoff" now shows "In Translation", and the (now narrower)
# Status columns that held that text are re-sized to fit the new content.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsZhCn     = $wb.Sheets.Item("zh-cn")
$wsDeDe     = $wb.Sheets.Item("de-de")

# Update every cell currently holding the old status text.
$wsOverview.Range("E2:F4").Value = "In Translation"
$wsZhCn.Range("C2:C4").Value     = "In Translation"
$wsDeDe.Range("C2:C4").Value     = "In Translation"

# The Status columns shrink to fit the new, shorter text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5   # column E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = 12.5   # column F (de-de status)
$wsZhCn.Columns.Item(3).ColumnWidth     = 12.5   # column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth     = 12.5   # column C (Status)
